# Update the "scraped_at" timestamps (column K) on the "snapshot" sheet
# to reflect the latest scrape run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$timestamps = @{
    2  = "2025-12-15T03:01:19.836255+00:00"
    3  = "2025-12-15T03:01:19.836294+00:00"
    4  = "2025-12-15T03:01:19.836316+00:00"
    5  = "2025-12-15T03:01:22.232285+00:00"
    6  = "2025-12-15T03:01:22.232314+00:00"
    7  = "2025-12-15T03:01:22.232332+00:00"
    8  = "2025-12-15T03:01:24.955934+00:00"
    9  = "2025-12-15T03:01:27.707940+00:00"
    10 = "2025-12-15T03:01:30.598338+00:00"
    11 = "2025-12-15T03:01:33.380600+00:00"
    12 = "2025-12-15T03:01:38.537032+00:00"
    13 = "2025-12-15T03:01:38.537061+00:00"
    14 = "2025-12-15T03:01:41.302659+00:00"
    15 = "2025-12-15T03:01:44.052153+00:00"
    16 = "2025-12-15T03:01:46.815694+00:00"
    17 = "2025-12-15T03:01:49.077241+00:00"
    18 = "2025-12-15T03:01:49.077269+00:00"
    19 = "2025-12-15T03:01:49.077287+00:00"
    20 = "2025-12-15T03:01:49.077305+00:00"
    21 = "2025-12-15T03:01:49.077321+00:00"
    22 = "2025-12-15T03:01:51.342829+00:00"
    23 = "2025-12-15T03:01:51.342857+00:00"
    24 = "2025-12-15T03:01:51.342875+00:00"
    25 = "2025-12-15T03:01:54.157468+00:00"
    26 = "2025-12-15T03:01:54.157501+00:00"
    27 = "2025-12-15T03:01:54.157522+00:00"
    28 = "2025-12-15T03:01:54.157552+00:00"
    29 = "2025-12-15T03:01:54.157569+00:00"
    30 = "2025-12-15T03:01:56.974508+00:00"
    31 = "2025-12-15T03:01:56.974539+00:00"
    32 = "2025-12-15T03:01:56.974559+00:00"
    33 = "2025-12-15T03:01:59.887471+00:00"
    34 = "2025-12-15T03:01:59.887502+00:00"
    35 = "2025-12-15T03:01:59.887520+00:00"
    36 = "2025-12-15T03:02:02.710171+00:00"
    37 = "2025-12-15T03:02:05.417754+00:00"
    38 = "2025-12-15T03:02:05.417783+00:00"
    39 = "2025-12-15T03:02:05.417802+00:00"
    40 = "2025-12-15T03:02:07.721622+00:00"
    41 = "2025-12-15T03:02:10.398658+00:00"
    42 = "2025-12-15T03:02:10.398685+00:00"
    43 = "2025-12-15T03:02:13.170298+00:00"
    44 = "2025-12-15T03:02:13.170333+00:00"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 11).Value = $timestamps[$row]
}
